# Doing Updates for Financials
# Update yearly figures on the ESTR financials worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESTR")

# Income Statement
$ws.Range("D8").Value  = 350000
$ws.Range("D9").Value  = 244500
$ws.Range("D10").Value = 105500
$ws.Range("D15").Value = 6700
$ws.Range("D17").Value = 320300
$ws.Range("D18").Value = 29700
$ws.Range("D20").Value = -49800
$ws.Range("D21").Value = 14300
$ws.Range("D22").Value = 59100
$ws.Range("D23").Value = -79100
$ws.Range("D24").Value = -90500
$ws.Range("D26").Value = 11300
$ws.Range("D27").Value = 9200
$ws.Range("D29").Value = 2100
$ws.Range("D32").Value = 49800
$ws.Range("D33").Value = 11200
$ws.Range("D35").Value = 11200

# Balance Sheet
$ws.Range("D41").Value = 11400
$ws.Range("D42").Value = 10400
$ws.Range("D43").Value = 206700
$ws.Range("D44").Value = 2900
$ws.Range("D45").Value = 1700
$ws.Range("D46").Value = 233000
$ws.Range("E46").Value = 300
$ws.Range("D47").Value = 50600
$ws.Range("E47").Value = 95000
$ws.Range("F47").Value = 94900
$ws.Range("D48").Value = 176800
$ws.Range("D49").Value = 150800
$ws.Range("D54").Value = 611200
$ws.Range("E54").Value = 95300
$ws.Range("F54").Value = 95300
$ws.Range("D57").Value = 32300
$ws.Range("D58").Value = 3600
$ws.Range("D59").Value = 121800
$ws.Range("D60").Value = 157800
$ws.Range("D61").Value = 369300
$ws.Range("D62").Value = 201000
$ws.Range("E62").Value = 3300
$ws.Range("F62").Value = 3300
$ws.Range("D66").Value = 726700
$ws.Range("E66").Value = 3400
$ws.Range("F66").Value = 3400
$ws.Range("D72").Value = -116000
$ws.Range("D76").Value = -115600
$ws.Range("E76").Value = 91900
$ws.Range("F76").Value = 91900

# Cash Flow Statement
$ws.Range("D81").Value = 11200
$ws.Range("D83").Value = 35100
$ws.Range("D89").Value = 62400
$ws.Range("D91").Value = -37200
$ws.Range("D94").Value = -51400
$ws.Range("F94").Value = -94900
$ws.Range("D100").Value = 2700
$ws.Range("F100").Value = 95300
$ws.Range("D102").Value = 13700
